# Natmi following Dr Hou advice:
# recomputed ligand-receptor pair statistics for FAPs<->ECs/FAPs/sCs and added the
# three new rows for the sCs sending cluster (Angpt4-Tek against ECs/FAPs/sCs).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colOrder = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$rowsData = @(
    @{ A="FAPs"; B="Angpt4"; C="Tek"; D="ECs"; E=3; F=1; G=1.741645; H=5.224935; I=0.903829134901074; J=0.9038291349010741; K=3; L=1; M=44.69746666666666; N=134.0924; O=0.6823972194925493; P=0.6823972194925493; Q=77.84711933266667; R=700.624073994; S=0.6167704885528491; T=0.6167704885528492 },
    @{ A="FAPs"; B="Angpt4"; C="Tek"; D="FAPs"; E=3; F=1; G=1.741645; H=5.224935; I=0.903829134901074; J=0.9038291349010741; K=3; L=1; M=18.63243533333333; N=55.897306; O=0.2844618053784121; P=0.2844618053784121; Q=32.45108783612334; R=292.05979052511; S=0.2571048674675679; T=0.2571048674675679 },
    @{ A="FAPs"; B="Angpt4"; C="Tek"; D="sCs"; E=3; F=1; G=1.741645; H=5.224935; I=0.903829134901074; J=0.9038291349010741; K=3; L=1; M=2.170755666666667; N=6.512267; O=0.03314097512903853; P=0.03314097512903853; Q=3.780685753071667; R=34.026171777645; S=0.02995377888065691; T=0.0299537788806569 },
    @{ A="sCs"; B="Angpt4"; C="Tek"; D="ECs"; E=2; F=0.6666666666666666; G=0.1853176666666667; H=0.555953; I=0.09617086509892599; J=0.09617086509892599; K=3; L=1; M=44.69746666666666; N=134.0924; O=0.6823972194925493; P=0.6823972194925493; Q=8.283230228577779; R=74.54907205720001; S=0.06562673093970015; T=0.06562673093970015 },
    @{ A="sCs"; B="Angpt4"; C="Tek"; D="FAPs"; E=2; F=0.6666666666666666; G=0.1853176666666667; H=0.555953; I=0.09617086509892599; J=0.09617086509892599; K=3; L=1; M=18.63243533333333; N=55.897306; O=0.2844618053784121; P=0.2844618053784121; Q=3.452919440290889; R=31.076274962618; S=0.02735693791084421; T=0.02735693791084421 },
    @{ A="sCs"; B="Angpt4"; C="Tek"; D="sCs"; E=2; F=0.6666666666666666; G=0.1853176666666667; H=0.555953; I=0.09617086509892599; J=0.09617086509892599; K=3; L=1; M=2.170755666666667; N=6.512267; O=0.03314097512903853; P=0.03314097512903853; Q=0.4022793750501112; R=3.620514375451001; S=0.003187196248381626; T=0.003187196248381625 },
)

$startRow = 2
for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $r = $startRow + $i
    $rowHash = $rowsData[$i]
    $colIndex = 0
    foreach ($col in $colOrder) {
        $colIndex = $colIndex + 1
        $ws.Cells.Item($r, $colIndex).Value = $rowHash[$col]
    }
}